$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update task text labels -------------------------------------------------
$ws.Range("C3").Value = "Phase 2 - Sprint4"
$ws.Range("C5").Value = "Implement -  Add constraint to task"
$ws.Range("C6").Value = "Implement - Add constraint date to task"
$ws.Range("C7").Value = "Implement - Show daily information system"
$ws.Range("C8").Value = "Implement - Close daily information system"

# --- Give C10 (currently blank) an underlined font, matching the row style --
$ws.Range("C10").Font.Underline = $true

# --- Remove the trailing blank rows, keep a lone formatted cell at D17 ------
$ws.Range("B15").Clear()
$ws.Rows("16").Delete()
$ws.Range("B16:F16").Clear()
$ws.Range("D16").Font.Underline = $true

# --- View / selection state ---------------------------------------------------
$ws.Range("C8").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
